$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.59454423119314
$ws.Cells.Item(2, 3).Value = 0.2002520751363761
$ws.Cells.Item(2, 4).Value = 0.4962892435936936
$ws.Cells.Item(2, 5).Value = 0.1607066697175679
$ws.Cells.Item(2, 7).Value = 0.8013244342879347
$ws.Cells.Item(2, 8).Value = 0.8868644463710922
$ws.Cells.Item(2, 9).Value = 0.8267295841770377
$ws.Cells.Item(2, 10).Value = 0.06435962793137051
$ws.Cells.Item(2, 12).Value = 0.4601634685019036
$ws.Cells.Item(2, 15).Value = 3.394152705191658
# Row 3
$ws.Cells.Item(3, 2).Value = 1.456537482687565
$ws.Cells.Item(3, 3).Value = 0.1829435412272744
$ws.Cells.Item(3, 4).Value = 0.4960566728477431
$ws.Cells.Item(3, 5).Value = 0.1619494273977473
$ws.Cells.Item(3, 7).Value = 0.8086428600207043
$ws.Cells.Item(3, 8).Value = 0.8952821707344256
$ws.Cells.Item(3, 9).Value = 0.8397790004629151
$ws.Cells.Item(3, 10).Value = 0.0646551489122924
$ws.Cells.Item(3, 12).Value = 0.4493776750452128
$ws.Cells.Item(3, 15).Value = 3.426965657155577
# Row 4
$ws.Cells.Item(4, 2).Value = 1.37173391705204
$ws.Cells.Item(4, 3).Value = 0.1722703648450761
$ws.Cells.Item(4, 4).Value = 0.4961431168577661
$ws.Cells.Item(4, 5).Value = 0.1627783241024012
$ws.Cells.Item(4, 7).Value = 0.8137551243897789
$ws.Cells.Item(4, 8).Value = 0.9009075401268944
$ws.Cells.Item(4, 9).Value = 0.8483320162841572
$ws.Cells.Item(4, 10).Value = 0.06485031052647372
$ws.Cells.Item(4, 12).Value = 0.4429033895283254
$ws.Cells.Item(4, 15).Value = 3.449369706683456
# Row 5
$ws.Cells.Item(5, 2).Value = 1.337161429472189
$ws.Cells.Item(5, 3).Value = 0.167909793390379
$ws.Cells.Item(5, 4).Value = 0.4962360746216348
$ws.Cells.Item(5, 5).Value = 0.1631326789052387
$ws.Cells.Item(5, 7).Value = 0.8159937855378132
$ws.Cells.Item(5, 8).Value = 0.9033148102327146
$ws.Cells.Item(5, 9).Value = 0.8519532596537545
$ws.Cells.Item(5, 10).Value = 0.06493329817005478
$ws.Cells.Item(5, 12).Value = 0.4403025720704932
$ws.Cells.Item(5, 15).Value = 3.459066581365775
# Row 6
$ws.Cells.Item(6, 2).Value = 1.331419900795822
$ws.Cells.Item(6, 3).Value = 0.167185058869336
$ws.Cells.Item(6, 4).Value = 0.4962549995388343
$ws.Cells.Item(6, 5).Value = 0.163192520699587
$ws.Cells.Item(6, 7).Value = 0.8163748914096942
$ws.Cells.Item(6, 8).Value = 0.9037214746680391
$ws.Cells.Item(6, 9).Value = 0.8525627639181987
$ws.Cells.Item(6, 10).Value = 0.06494728733896959
$ws.Cells.Item(6, 12).Value = 0.4398729804483281
$ws.Cells.Item(6, 15).Value = 3.460710975103282
# Row 7
$ws.Cells.Item(7, 2).Value = 1.37126771519678
$ws.Cells.Item(7, 3).Value = 0.1722116014380788
$ws.Cells.Item(7, 4).Value = 0.4961441366531147
$ws.Cells.Item(7, 5).Value = 0.1627830359303424
$ws.Cells.Item(7, 7).Value = 0.813784686943066
$ws.Cells.Item(7, 8).Value = 0.9009395402352354
$ws.Cells.Item(7, 9).Value = 0.8483803038741176
$ws.Cells.Item(7, 10).Value = 0.06485141571357822
$ws.Cells.Item(7, 12).Value = 0.4428681618237249
$ws.Cells.Item(7, 15).Value = 3.449498186781526
# Row 8
$ws.Cells.Item(8, 2).Value = 1.546974813536792
$ws.Cells.Item(8, 3).Value = 0.1942937553513104
$ws.Cells.Item(8, 4).Value = 0.496161524037035
$ws.Cells.Item(8, 5).Value = 0.1611215180402024
$ws.Cells.Item(8, 7).Value = 0.8037192848971486
$ws.Cells.Item(8, 8).Value = 0.8896720623748706
$ws.Cells.Item(8, 9).Value = 0.8311167439082361
$ws.Cells.Item(8, 10).Value = 0.06445868397692323
$ws.Cells.Item(8, 12).Value = 0.4564138950918846
$ws.Cells.Item(8, 15).Value = 3.404997955104321
# Row 9
$ws.Cells.Item(9, 2).Value = 1.89091444080691
$ws.Cells.Item(9, 3).Value = 0.2372226077008861
$ws.Cells.Item(9, 4).Value = 0.4980113562771891
$ws.Cells.Item(9, 5).Value = 0.1583849395572461
$ws.Cells.Item(9, 7).Value = 0.7889002131679774
$ws.Cells.Item(9, 8).Value = 0.8712012595851064
$ws.Cells.Item(9, 9).Value = 0.8015577178457498
$ws.Cells.Item(9, 10).Value = 0.06379688086221513
$ws.Cells.Item(9, 12).Value = 0.4841447559719825
$ws.Cells.Item(9, 15).Value = 3.335660337338084
# Row 10
$ws.Cells.Item(10, 2).Value = 2.143131559468145
$ws.Cells.Item(10, 3).Value = 0.2685215305963595
$ws.Cells.Item(10, 4).Value = 0.5004738819357328
$ws.Cells.Item(10, 5).Value = 0.1566913919459161
$ws.Cells.Item(10, 7).Value = 0.7810256979376931
$ws.Cells.Item(10, 8).Value = 0.8598402759256913
$ws.Cells.Item(10, 9).Value = 0.7824654851299861
$ws.Cells.Item(10, 10).Value = 0.06337610262314897
$ws.Cells.Item(10, 12).Value = 0.5052216395487079
$ws.Cells.Item(10, 15).Value = 3.295677097891087
# Row 11
$ws.Cells.Item(11, 2).Value = 2.257749858687419
$ws.Cells.Item(11, 3).Value = 0.2827054323715572
$ws.Cells.Item(11, 4).Value = 0.5018330762063385
$ws.Cells.Item(11, 5).Value = 0.1559895816481234
$ws.Cells.Item(11, 7).Value = 0.7781007153527071
$ws.Cells.Item(11, 8).Value = 0.8551516306786766
$ws.Cells.Item(11, 9).Value = 0.7743512763347056
$ws.Cells.Item(11, 10).Value = 0.06319876545307679
$ws.Cells.Item(11, 12).Value = 0.5149610106036135
$ws.Cells.Item(11, 15).Value = 3.279873957306961
# Row 12
$ws.Cells.Item(12, 2).Value = 2.301133980230418
$ws.Cells.Item(12, 3).Value = 0.2880684513825713
$ws.Cells.Item(12, 4).Value = 0.5023820604756679
$ws.Cells.Item(12, 5).Value = 0.1557336716967672
$ws.Cells.Item(12, 7).Value = 0.7770878432330619
$ws.Cells.Item(12, 8).Value = 0.8534451218849455
$ws.Cells.Item(12, 9).Value = 0.771360872013986
$ws.Cells.Item(12, 10).Value = 0.063133626843209
$ws.Cells.Item(12, 12).Value = 0.5186706320100853
$ws.Cells.Item(12, 15).Value = 3.274233271398714
# Row 13
$ws.Cells.Item(13, 2).Value = 2.291791334505774
$ws.Cells.Item(13, 3).Value = 0.2869137955317456
$ws.Cells.Item(13, 4).Value = 0.5022623032136693
$ws.Cells.Item(13, 5).Value = 0.1557883485954701
$ws.Cells.Item(13, 7).Value = 0.7773017650631147
$ws.Cells.Item(13, 8).Value = 0.8538095808181367
$ws.Cells.Item(13, 9).Value = 0.7720012473940585
$ws.Cells.Item(13, 10).Value = 0.06314756612915495
$ws.Cells.Item(13, 12).Value = 0.5178707447411313
$ws.Cells.Item(13, 15).Value = 3.275432801939019
# Row 14
$ws.Cells.Item(14, 2).Value = 2.261319498816704
$ws.Cells.Item(14, 3).Value = 0.2831468157357335
$ws.Cells.Item(14, 4).Value = 0.501877554759119
$ws.Cells.Item(14, 5).Value = 0.1559683304417092
$ws.Cells.Item(14, 7).Value = 0.7780154852643335
$ws.Cells.Item(14, 8).Value = 0.8550098525822563
$ws.Cells.Item(14, 9).Value = 0.7741036041796079
$ws.Cells.Item(14, 10).Value = 0.06319336612583193
$ws.Cells.Item(14, 12).Value = 0.5152657731888581
$ws.Cells.Item(14, 15).Value = 3.279403004604916
# Row 15
$ws.Cells.Item(15, 2).Value = 2.242652024551603
$ws.Cells.Item(15, 3).Value = 0.2808383649970949
$ws.Cells.Item(15, 4).Value = 0.5016463480720574
$ws.Cells.Item(15, 5).Value = 0.156079856937744
$ws.Cells.Item(15, 7).Value = 0.7784650069597632
$ws.Cells.Item(15, 8).Value = 0.8557540375788335
$ws.Cells.Item(15, 9).Value = 0.7754020774803152
$ws.Cells.Item(15, 10).Value = 0.06322168210671819
$ws.Cells.Item(15, 12).Value = 0.513672950114028
$ws.Cells.Item(15, 15).Value = 3.281879635584204
# Row 16
$ws.Cells.Item(16, 2).Value = 2.135638417764881
$ws.Cells.Item(16, 3).Value = 0.2675934614484845
$ws.Cells.Item(16, 4).Value = 0.5003898575943992
$ws.Cells.Item(16, 5).Value = 0.1567386360356515
$ws.Cells.Item(16, 7).Value = 0.7812300982386375
$ws.Cells.Item(16, 8).Value = 0.8601563412505584
$ws.Cells.Item(16, 9).Value = 0.7830072710943909
$ws.Cells.Item(16, 10).Value = 0.06338797450696632
$ws.Cells.Item(16, 12).Value = 0.5045881748422971
$ws.Cells.Item(16, 15).Value = 3.296757922528116
# Row 17
$ws.Cells.Item(17, 2).Value = 2.069957344060583
$ws.Cells.Item(17, 3).Value = 0.259454039408098
$ws.Cells.Item(17, 4).Value = 0.4996801941410638
$ws.Cells.Item(17, 5).Value = 0.1571603345749839
$ws.Cells.Item(17, 7).Value = 0.7830948917588501
$ws.Cells.Item(17, 8).Value = 0.8629798468234213
$ws.Cells.Item(17, 9).Value = 0.7878191640600818
$ws.Cells.Item(17, 10).Value = 0.06349358843460884
$ws.Cells.Item(17, 12).Value = 0.4990535653969914
$ws.Cells.Item(17, 15).Value = 3.306496663305097
# Row 18
$ws.Cells.Item(18, 2).Value = 2.032168512044962
$ws.Cells.Item(18, 3).Value = 0.2547673833179545
$ws.Cells.Item(18, 4).Value = 0.4992945129920656
$ws.Cells.Item(18, 5).Value = 0.1574093415077531
$ws.Cells.Item(18, 7).Value = 0.7842293128688453
$ws.Cells.Item(18, 8).Value = 0.864648986069497
$ws.Cells.Item(18, 9).Value = 0.7906405674235657
$ws.Cells.Item(18, 10).Value = 0.06355566030163828
$ws.Cells.Item(18, 12).Value = 0.4958844681513312
$ws.Cells.Item(18, 15).Value = 3.312322622369976
# Row 19
$ws.Cells.Item(19, 2).Value = 2.019372084315876
$ws.Cells.Item(19, 3).Value = 0.2531797004918133
$ws.Cells.Item(19, 4).Value = 0.4991677944964863
$ws.Cells.Item(19, 5).Value = 0.157494760390632
$ws.Cells.Item(19, 7).Value = 0.7846240218720482
$ws.Cells.Item(19, 8).Value = 0.8652218789184474
$ws.Cells.Item(19, 9).Value = 0.7916050691036887
$ws.Cells.Item(19, 10).Value = 0.06357690470989574
$ws.Cells.Item(19, 12).Value = 0.4948139233220274
$ws.Cells.Item(19, 15).Value = 3.314333730628249
# Row 20
$ws.Cells.Item(20, 2).Value = 2.076950344642455
$ws.Cells.Item(20, 3).Value = 0.2603210225246926
$ws.Cells.Item(20, 4).Value = 0.4997534110828781
$ws.Cells.Item(20, 5).Value = 0.1571147758573748
$ws.Cells.Item(20, 7).Value = 0.7828899791233823
$ws.Cells.Item(20, 8).Value = 0.8626746085058983
$ws.Cells.Item(20, 9).Value = 0.7873013677768768
$ws.Cells.Item(20, 10).Value = 0.06348220853050179
$ws.Cells.Item(20, 12).Value = 0.4996412594996684
$ws.Cells.Item(20, 15).Value = 3.305436719199804
# Row 21
$ws.Cells.Item(21, 2).Value = 2.270270360992868
$ws.Cells.Item(21, 3).Value = 0.2842534919327875
$ws.Cells.Item(21, 4).Value = 0.5019896347915562
$ws.Cells.Item(21, 5).Value = 0.1559151981794038
$ws.Cells.Item(21, 7).Value = 0.7778032746883383
$ws.Cells.Item(21, 8).Value = 0.854655431393482
$ws.Cells.Item(21, 9).Value = 0.773483856884372
$ws.Cells.Item(21, 10).Value = 0.06317985893912592
$ws.Cells.Item(21, 12).Value = 0.5160303339629877
$ws.Cells.Item(21, 15).Value = 3.278227529335538
# Row 22
$ws.Cells.Item(22, 2).Value = 2.396502250463072
$ws.Cells.Item(22, 3).Value = 0.2998472816076401
$ws.Cells.Item(22, 4).Value = 0.5036509503564588
$ws.Cells.Item(22, 5).Value = 0.1551886142516405
$ws.Cells.Item(22, 7).Value = 0.7750312427201465
$ws.Cells.Item(22, 8).Value = 0.8498165117506886
$ws.Cells.Item(22, 9).Value = 0.7649329283364885
$ws.Cells.Item(22, 10).Value = 0.06299399754545476
$ws.Cells.Item(22, 12).Value = 0.5268669430553103
$ws.Cells.Item(22, 15).Value = 3.262447908896092
# Row 23
$ws.Cells.Item(23, 2).Value = 2.329141186723575
$ws.Cells.Item(23, 3).Value = 0.2915290343340757
$ws.Cells.Item(23, 4).Value = 0.5027460177339407
$ws.Cells.Item(23, 5).Value = 0.1555711571576666
$ws.Cells.Item(23, 7).Value = 0.7764600978872522
$ws.Cells.Item(23, 8).Value = 0.8523623361189294
$ws.Cells.Item(23, 9).Value = 0.7694527792394403
$ws.Cells.Item(23, 10).Value = 0.06309212390603136
$ws.Cells.Item(23, 12).Value = 0.5210718445394065
$ws.Cells.Item(23, 15).Value = 3.27068630357337
# Row 24
$ws.Cells.Item(24, 2).Value = 2.073788895813323
$ws.Cells.Item(24, 3).Value = 0.259929081842273
$ws.Cells.Item(24, 4).Value = 0.4997202401932981
$ws.Cells.Item(24, 5).Value = 0.1571353524977894
$ws.Cells.Item(24, 7).Value = 0.7829824260012828
$ws.Cells.Item(24, 8).Value = 0.8628124638721459
$ws.Cells.Item(24, 9).Value = 0.7875352922221026
$ws.Cells.Item(24, 10).Value = 0.06348734917001408
$ws.Cells.Item(24, 12).Value = 0.4993755230357095
$ws.Cells.Item(24, 15).Value = 3.305915212793394
# Row 25
$ws.Cells.Item(25, 2).Value = 1.797946249039455
$ws.Cells.Item(25, 3).Value = 0.2256506097980377
$ws.Cells.Item(25, 4).Value = 0.4973168993132617
$ws.Cells.Item(25, 5).Value = 0.1559683304417092
$ws.Cells.Item(25, 7).Value = 0.792381061626088
$ws.Cells.Item(25, 8).Value = 0.8758100735997942
$ws.Cells.Item(25, 9).Value = 0.8090937528665592
$ws.Cells.Item(25, 10).Value = 0.06396437971413249
$ws.Cells.Item(25, 12).Value = 0.4765187311326855
$ws.Cells.Item(25, 15).Value = 3.352495637072678
